# Refresh the daily crypto snapshot (Price / Volume(1h) columns, plus the
# #49 row which rotated from Hedera to InjectiveProtocol) to match the
# latest coinranking.com pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '60.990.84'
$ws.Range("E2").Value = '  -1.51%  '

# Row 3
$ws.Range("D3").Value = '2.431.50'
$ws.Range("E3").Value = '  -0.09%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("E4").Value = '  -0.10%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '571.84'
$ws.Range("E5").Value = '  -1.45%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '140.43'
$ws.Range("E6").Value = '  -1.53%  '

# Row 7
$ws.Range("E7").Value = '  +0.18%  '

# Row 8
$ws.Range("E8").Value = '  +0.22%  '

# Row 9
$ws.Range("D9").Value = '2.419.78'
$ws.Range("E9").Value = '  -0.43%  '

# Row 10
$ws.Range("E10").Value = '  +1.88%  '

# Row 11
$ws.Range("E11").Value = '  +1.15%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.12'
$ws.Range("E12").Value = '  -1.39%  '

# Row 13
$ws.Range("E13").Value = '  -1.06%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.09'
$ws.Range("E14").Value = '  -0.66%  '

# Row 15
$ws.Range("D15").Value = '2.869.19'
$ws.Range("E15").Value = '  -0.19%  '

# Row 16
$ws.Range("E16").Value = '  -0.62%  '

# Row 17
$ws.Range("D17").Value = '61.126.18'
$ws.Range("E17").Value = '  -1.23%  '

# Row 18
$ws.Range("D18").Value = '2.449.10'
$ws.Range("E18").Value = '  +0.71%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.58'
$ws.Range("E19").Value = '  -2.51%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.26'
$ws.Range("E20").Value = '  +2.76%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '323.56'
$ws.Range("E21").Value = '  -1.85%  '

# Row 22
$ws.Range("E22").Value = '  -1.12%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.10'
$ws.Range("E23").Value = '  +2.25%  '

# Row 24
$ws.Range("E24").Value = '  +0.15%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.90'
$ws.Range("E25").Value = '  -2.37%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '65.14'
$ws.Range("E26").Value = '  -0.64%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.87'
$ws.Range("E27").Value = '  -4.26%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '573.57'
$ws.Range("E28").Value = '  -6.42%  '

# Row 29
$ws.Range("D29").Value = '2.549.80'
$ws.Range("E29").Value = '  -0.17%  '

# Row 30
$ws.Range("E30").Value = '  -0.05%  '

# Row 31
$ws.Range("D31").Value = '0.0₃0909'
$ws.Range("E31").Value = '  -3.85%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.88'
$ws.Range("E32").Value = '  -1.01%  '

# Row 33
$ws.Range("E33").Value = '  -5.52%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.84'
$ws.Range("E34").Value = '  -1.89%  '

# Row 35
$ws.Range("E35").Value = '  -6.32%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.01'
$ws.Range("E36").Value = '  +0.37%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.61'
$ws.Range("E37").Value = '  -5.59%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '152.09'
$ws.Range("E38").Value = '  +0.58%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.368'
$ws.Range("E39").Value = '  -1.62%  '

# Row 40
$ws.Range("E40").Value = '  -2.78%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '18.29'
$ws.Range("E41").Value = '  +0.00%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.10'
$ws.Range("E42").Value = '  -1.99%  '

# Row 43
$ws.Range("E43").Value = '  +0.00%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '41.70'

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.65'
$ws.Range("E45").Value = '  -5.21%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.34'
$ws.Range("E46").Value = '  -4.04%  '

# Row 47
$ws.Range("D47").Value = '0.0₆0286'
$ws.Range("E47").Value = '  +26.05%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '141.01'
$ws.Range("E48").Value = '  -1.10%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.52'
$ws.Range("E49").Value = '  -2.23%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.594'
$ws.Range("E50").Value = '  -0.21%  '

# Row 51
$ws.Range("B51").Value = 'InjectiveProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '19.49'
$ws.Range("E51").Value = '  +0.48%  '
